# "add customer ui and guest ui"
# The author filled in passwords for the two customer rows on the
# "Customer" sheet (Nancy / Jerry) and left that sheet active/selected
# instead of "Room".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer")

# Nancy's password 123 -> 123456
$ws.Range("F2").Value = 123456
# Jerry's password was blank -> 123456
$ws.Range("F3").Value = 123456

# Make "Customer" the active/selected sheet (it was "Room" before) and
# leave the cursor on E19, matching the saved selection.
$ws.Activate() | Out-Null
$ws.Range("E19").Select() | Out-Null
